$d = $word.ActiveDocument

# Remove both reviewer comments (and, as a side effect, the
# commentRangeStart / commentRangeEnd / commentReference runs that
# anchor them in the body). Deleting index 1 twice removes comment 0
# (Samantha Jane Bolten - "Please Note..." paragraph) followed by what
# was comment 4 (Caroline Motzer - "No signs of improvement"). This
# also renumbers the surviving bookmarks (_Hlk124086675/656/628) down
# by one id, same as Word does when the comment id "slot" in front of
# them disappears.
while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# The second comment used to split "No control/No signs of
# improvement" into two separate runs around the comment range. Now
# that the comment markers are gone, collapse the wording back into a
# single run of text. Scope the Find/Replace to that one list item
# (there's a second, unrelated "No control/No signs of improvement"
# bullet later in the survey that must stay untouched), and do it with
# revision tracking briefly suspended so the merge lands as a plain
# edit instead of an w:ins/w:del pair - the document's TrackRevisions
# setting itself is left exactly as the author had it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "No control/No signs of improvement*") {
        $target = $p.Range
        break
    }
}
if ($target -ne $null) {
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $false
    $target.Find.Execute("No control/No signs of improvement", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "No control/No signs of improvement", 2) | Out-Null
    $d.TrackRevisions = $wasTracking
}
